# Update "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet (which aggregates 展览 + 演出 rows, so its row numbers
# are shifted by +1 starting at row 23 because of one extra 演出 row).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for "展览"
$exhibitUpdates = @{
    2  = 593
    5  = 1118
    6  = 14150
    7  = 15810
    17 = 32
    20 = 1225
    23 = 19
    24 = 6243
    27 = 5594
    31 = 4583
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Row -> new F value for "全部类型" (rows offset by +1 vs 展览 from row 23 on)
$allUpdates = @{
    2  = 593
    5  = 1118
    6  = 14150
    7  = 15810
    17 = 32
    20 = 1225
    24 = 19
    25 = 6243
    28 = 5594
    32 = 4583
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
